# Fix systematic spacing issue between header bar and body text
# (content-level part of the change: condense CORE COMPETENCIES into a
# single summary line, and re-add the detailed category breakdowns as a
# new TECHNICAL SKILLS section at the end of the document.)

$d = $word.ActiveDocument

$bullet = [char]0x2022

# ---------------------------------------------------------------------
# 1. Collapse the three detailed "CORE COMPETENCIES" paragraphs into a
#    single condensed paragraph that just lists the three category
#    headings separated by bullets.
# ---------------------------------------------------------------------

$coreHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "CORE COMPETENCIES") {
        $coreHeading = $p
        break
    }
}

$firstCompetency = $coreHeading.Next()
$secondCompetency = $firstCompetency.Next()
$thirdCompetency = $secondCompetency.Next()

# Remove the 2nd and 3rd detail paragraphs entirely (including their
# paragraph marks), then rewrite the 1st into the condensed summary.
$thirdCompetency.Range.Delete()
$secondCompetency.Range.Delete()

$firstCompetency.Range.Text = "Statistical Analysis & Machine Learning $bullet Big Data & Data Engineering $bullet Data Visualization & Reporting"

# ---------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the very end of the
#    document body, restoring the detailed category breakdowns in a
#    condensed "Category; Category; Category" style.
# ---------------------------------------------------------------------

$endRange = $d.Paragraphs.Last.Range
$endRange.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading 2"

$headingPara.Range.InsertParagraphAfter()
$statPara = $d.Paragraphs.Last
$statPara.Style = "Normal"
$statPara.Range.Text = "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques"

$statPara.Range.InsertParagraphAfter()
$bigDataPara = $d.Paragraphs.Last
$bigDataPara.Style = "Normal"
$bigDataPara.Range.Text = "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization"

$bigDataPara.Range.InsertParagraphAfter()
$dataVizPara = $d.Paragraphs.Last
$dataVizPara.Style = "Normal"
$dataVizPara.Range.Text = "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation"
